$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 with same style as the other header cells (e.g. E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$timestamps = @(
    "2021-10-05 10:51:33.208978",
    "2021-10-05 10:51:33.208991",
    "2021-10-05 10:51:33.208995",
    "2021-10-05 10:51:33.208999",
    "2021-10-05 10:51:33.209002",
    "2021-10-05 10:51:33.209006",
    "2021-10-05 10:51:33.209009",
    "2021-10-05 10:51:33.209012",
    "2021-10-05 10:51:33.209015",
    "2021-10-05 10:51:33.209018",
    "2021-10-05 10:51:33.209021",
    "2021-10-05 10:51:33.209025",
    "2021-10-05 10:51:33.209028",
    "2021-10-05 10:51:33.209031",
    "2021-10-05 10:51:33.209034",
    "2021-10-05 10:51:33.209037",
    "2021-10-05 10:51:33.209040",
    "2021-10-05 10:51:33.209044",
    "2021-10-05 10:51:33.209047",
    "2021-10-05 10:51:33.209050",
    "2021-10-05 10:51:33.209053",
    "2021-10-05 10:51:33.209056",
    "2021-10-05 10:51:33.209059",
    "2021-10-05 10:51:33.209063",
    "2021-10-05 10:51:33.209066",
    "2021-10-05 10:51:33.209069",
    "2021-10-05 10:51:33.209072",
    "2021-10-05 10:51:33.209075",
    "2021-10-05 10:51:33.209078",
    "2021-10-05 10:51:33.209081",
    "2021-10-05 10:51:33.209085",
    "2021-10-05 10:51:33.209088",
    "2021-10-05 10:51:33.209091",
    "2021-10-05 10:51:33.209094",
    "2021-10-05 10:51:33.209097",
    "2021-10-05 10:51:33.209100",
    "2021-10-05 10:51:33.209103",
    "2021-10-05 10:51:33.209106",
    "2021-10-05 10:51:33.209110",
    "2021-10-05 10:51:33.209113",
    "2021-10-05 10:51:33.209116",
    "2021-10-05 10:51:33.209120",
    "2021-10-05 10:51:33.209123",
    "2021-10-05 10:51:33.209126",
    "2021-10-05 10:51:33.209129",
    "2021-10-05 10:51:33.209132",
    "2021-10-05 10:51:33.209135",
    "2021-10-05 10:51:33.209138",
    "2021-10-05 10:51:33.209141",
    "2021-10-05 10:51:33.209144",
    "2021-10-05 10:51:33.209147",
    "2021-10-05 10:51:33.209150",
    "2021-10-05 10:51:33.209154",
    "2021-10-05 10:51:33.209157",
    "2021-10-05 10:51:33.209160",
    "2021-10-05 10:51:33.209163",
    "2021-10-05 10:51:33.209166"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
